$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all data cells are treated as text so date-like and numeric-like
# strings (e.g. "10/10/2023", "2023-11-11") are preserved literally instead
# of being auto-converted by Excel into date/number values.
$ws.Range("A2:E18").NumberFormat = "@"

# Row 2: Sarah Adams
$ws.Cells.Item(2, 1).Value = "Sarah Adams"
$ws.Cells.Item(2, 2).Value = "10/10/2023"
$ws.Cells.Item(2, 3).Value = "96 bpm"
$ws.Cells.Item(2, 4).Value = "140/70 mmHg"
$ws.Cells.Item(2, 5).Value = "3 hours"

# Row 3: Michael Brown
$ws.Cells.Item(3, 1).Value = "Michael Brown"
$ws.Cells.Item(3, 2).Value = "11/09/2023"
$ws.Cells.Item(3, 3).Value = "110 bpm"
$ws.Cells.Item(3, 4).Value = "120/85 mmHg"
$ws.Cells.Item(3, 5).Value = "6.5 hours"

# Row 4: Emily White
$ws.Cells.Item(4, 1).Value = "Emily White"
$ws.Cells.Item(4, 2).Value = "09/09/2023"
$ws.Cells.Item(4, 3).Value = "90 bpm"
$ws.Cells.Item(4, 4).Value = "130/90 mmHg"
$ws.Cells.Item(4, 5).Value = "5 hours"

# Row 5: David Lee
$ws.Cells.Item(5, 1).Value = "David Lee"
$ws.Cells.Item(5, 2).Value = "10/10/2023"
$ws.Cells.Item(5, 3).Value = "60 bpm"
$ws.Cells.Item(5, 4).Value = "120/75 mmHg"
$ws.Cells.Item(5, 5).Value = "6.5 hours"

# Row 6: Emma Taylor
$ws.Cells.Item(6, 1).Value = "Emma Taylor"
$ws.Cells.Item(6, 2).Value = "05/11/2023"
$ws.Cells.Item(6, 3).Value = "75 bpm"
$ws.Cells.Item(6, 4).Value = "115/95 mmHg"
$ws.Cells.Item(6, 5).Value = "9 hours"

# Row 7: Oliver Smith
$ws.Cells.Item(7, 1).Value = "Oliver Smith"
$ws.Cells.Item(7, 2).Value = "03/09/2023"
$ws.Cells.Item(7, 3).Value = "120 bpm"
$ws.Cells.Item(7, 4).Value = "150/95 mmHg"
$ws.Cells.Item(7, 5).Value = "10.5 hours"

# Row 8: Sophia Miller
$ws.Cells.Item(8, 1).Value = "Sophia Miller"
$ws.Cells.Item(8, 2).Value = "26/10/2023"
$ws.Cells.Item(8, 3).Value = "80 bpm"
$ws.Cells.Item(8, 4).Value = "110/55 mmHg"
$ws.Cells.Item(8, 5).Value = "7 hours"

# Row 9: Sophia Wilson
$ws.Cells.Item(9, 1).Value = "Sophia Wilson"
$ws.Cells.Item(9, 2).Value = "2023-11-11"
$ws.Cells.Item(9, 3).Value = "72 bpm"
$ws.Cells.Item(9, 4).Value = "125/80 mmHg"
$ws.Cells.Item(9, 5).Value = "7 hours"

# Row 10: Liam Brown
$ws.Cells.Item(10, 1).Value = "Liam Brown"
$ws.Cells.Item(10, 2).Value = "2023-11-09"
$ws.Cells.Item(10, 3).Value = "82 bpm"
$ws.Cells.Item(10, 4).Value = "130/75 mmHg"
$ws.Cells.Item(10, 5).Value = "5.5 hours"

# Row 11: Olivia Taylor
$ws.Cells.Item(11, 1).Value = "Olivia Taylor"
$ws.Cells.Item(11, 2).Value = "2023-11-08"
$ws.Cells.Item(11, 3).Value = "62 bpm"
$ws.Cells.Item(11, 4).Value = "140/85 mmHg"
$ws.Cells.Item(11, 5).Value = "7.5 hours"

# Row 12: Noah Davis
$ws.Cells.Item(12, 1).Value = "Noah Davis"
$ws.Cells.Item(12, 2).Value = "2023-11-02"
$ws.Cells.Item(12, 3).Value = "122 bpm"
$ws.Cells.Item(12, 4).Value = "70/45 mmHg"
$ws.Cells.Item(12, 5).Value = "9 hours"

# Row 13: Ella Brown
$ws.Cells.Item(13, 1).Value = "Ella Brown"
$ws.Cells.Item(13, 2).Value = "2023-11-07"
$ws.Cells.Item(13, 3).Value = "67 bpm"
$ws.Cells.Item(13, 4).Value = "90/60 mmHg"
$ws.Cells.Item(13, 5).Value = "10 hours"

# Row 14: Charlotte Anderson
$ws.Cells.Item(14, 1).Value = "Charlotte Anderson"
$ws.Cells.Item(14, 2).Value = "2023-11-08"
$ws.Cells.Item(14, 3).Value = "87 bpm"
$ws.Cells.Item(14, 4).Value = "135/75 mmHg"
$ws.Cells.Item(14, 5).Value = "8 hours"

# Row 15: Ella Johnson
$ws.Cells.Item(15, 1).Value = "Ella Johnson"
$ws.Cells.Item(15, 2).Value = "2023-10-10"
$ws.Cells.Item(15, 3).Value = "96 bpm"
$ws.Cells.Item(15, 4).Value = "140/70 mmHg"
$ws.Cells.Item(15, 5).Value = "3 hours"

# Row 16: William Miller
$ws.Cells.Item(16, 1).Value = "William Miller"
$ws.Cells.Item(16, 2).Value = "2023-11-09"
$ws.Cells.Item(16, 3).Value = "86 bpm"
$ws.Cells.Item(16, 4).Value = "120/85 mmHg"
$ws.Cells.Item(16, 5).Value = "9 hours"

# Row 17: Scarlett Smith
$ws.Cells.Item(17, 1).Value = "Scarlett Smith"
$ws.Cells.Item(17, 2).Value = "2023-01-09"
$ws.Cells.Item(17, 3).Value = "100 bpm"
$ws.Cells.Item(17, 4).Value = "130/65 mmHg"
$ws.Cells.Item(17, 5).Value = "8 hours"

# Row 18: James Lee
$ws.Cells.Item(18, 1).Value = "James Lee"
$ws.Cells.Item(18, 2).Value = "2023-09-29"
$ws.Cells.Item(18, 3).Value = "125 bpm"
$ws.Cells.Item(18, 4).Value = "160/95 mmHg"
$ws.Cells.Item(18, 5).Value = "5.5 hours"
